# Adds a new section "Agregando la fecha de publicación" at the end of the
# document, right before the trailing (empty) paragraph that carries the
# document's "_GoBack" bookmark.
#
# New content, in order:
#   1. An empty paragraph (paragraph-mark only, underline formatting) that
#      mirrors the separator blank lines used elsewhere in this document.
#   2. A heading paragraph: "Agregando la fecha de publicación".
#   3. A body paragraph describing Moment.js / Platzigram relative dates.
#
# The trailing bookmark paragraph itself is left untouched (re-emitted
# verbatim) since Range.InsertXML on the single-character paragraph-mark
# range replaces that range's content.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---- shared run-properties blocks -----------------------------------
$sepRpr = '<w:u w:val="single"/>'

$headingRpr = '<w:rFonts w:ascii="cooper_hewittmedium" w:eastAsia="Times New Roman" w:hAnsi="cooper_hewittmedium" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="000000"/><w:spacing w:val="-2"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:eastAsia="es-CO"/>'

$bodyRpr = '<w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="273B47"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="es-CO"/>'

# ---- paragraph 1: blank separator line --------------------------------
$separatorPara = "<w:p $wNs><w:pPr><w:rPr>$sepRpr</w:rPr></w:pPr></w:p>"

# ---- paragraph 2: heading ----------------------------------------------
$headingText = 'Agregando la fecha de publicación'
$headingPara = "<w:p $wNs>" +
  "<w:pPr><w:spacing w:line=`"240`" w:lineRule=`"auto`"/><w:rPr>$headingRpr</w:rPr></w:pPr>" +
  "<w:r><w:rPr>$headingRpr</w:rPr><w:t>$headingText</w:t></w:r>" +
  "</w:p>"

# ---- paragraph 3: body text (with proofErr spans) ----------------------
$run1 = 'Moment.js es una librería que nos ofrece diferentes formas de agregar fechas, desde el tradicional formato 24 horas hasta las fechas relativas (Hace x tiempo). '
$run2 = 'Platzigram'
$run3 = ' usa fechas relativas, por lo tanto, solo necesitamos pasar como parámetro la fecha actual '
$run4 = '( se'
$run5 = ' puede usar new Date() ) y la librería se encarga del resto. '

$bodyPara = "<w:p $wNs>" +
  "<w:pPr><w:spacing w:after=`"0`" w:line=`"240`" w:lineRule=`"auto`"/><w:rPr>$bodyRpr</w:rPr></w:pPr>" +
  "<w:r><w:rPr>$bodyRpr</w:rPr><w:t xml:space=`"preserve`">$run1</w:t></w:r>" +
  '<w:proofErr w:type="spellStart"/>' +
  "<w:r><w:rPr>$bodyRpr</w:rPr><w:t>$run2</w:t></w:r>" +
  '<w:proofErr w:type="spellEnd"/>' +
  "<w:r><w:rPr>$bodyRpr</w:rPr><w:t xml:space=`"preserve`">$run3</w:t></w:r>" +
  '<w:proofErr w:type="gramStart"/>' +
  "<w:r><w:rPr>$bodyRpr</w:rPr><w:t>$run4</w:t></w:r>" +
  '<w:proofErr w:type="gramEnd"/>' +
  "<w:r><w:rPr>$bodyRpr</w:rPr><w:t xml:space=`"preserve`">$run5</w:t></w:r>" +
  "</w:p>"

# ---- re-emit the trailing bookmark paragraph verbatim ------------------
$finalPara = "<w:p $wNs w:rsidR=`"00C776B8`" w:rsidRPr=`"006E7FB9`" w:rsidRDefault=`"00C776B8`">" +
  "<w:pPr><w:rPr>$sepRpr</w:rPr></w:pPr>" +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  "</w:p>"

$newXml = $separatorPara + $headingPara + $bodyPara + $finalPara

# The very last paragraph in the document is the (empty) bookmark
# paragraph; its Range covers exactly its own paragraph mark. Replacing
# that range with our XML block inserts the new paragraphs right before
# it while re-creating it unchanged as the new last paragraph.
$lastParaIndex = $d.Paragraphs.Count
$anchor = $d.Paragraphs.Item($lastParaIndex).Range
$anchor.InsertXML($newXml)
